$d = $word.ActiveDocument

$replacements = @(
    @("57×41=", "68×82="),
    @("94×74=", "50×74="),
    @("81×54=", "64×42="),
    @("61×87=", "34×51="),
    @("79×41=", "51×91="),
    @("66×98=", "76×27="),
    @("40×40=", "83×72="),
    @("23×86=", "41×14="),
    @("97×83=", "17×35="),
    @("74×82=", "16×46="),
    @("31×37=", "21×41="),
    @("27×76=", "25×76="),
    @("61×75=", "87×69="),
    @("69×82=", "71×60="),
    @("85×86=", "81×75="),
    @("58×91=", "72×85="),
    @("45×17=", "36×50="),
    @("57×22=", "26×32="),
    @("28×89=", "49×67="),
    @("85×53=", "39×67="),
    @("24×70=", "53×69="),
    @("96×93=", "13×18="),
    @("67×44=", "55×52="),
    @("39×61=", "55×96="),
    @("20×12=", "56×89=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
